# revisions before resubmission - lagged analysis, etc
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TableS3_prey_stocks_all")

$ram = "RAM Legacy Database"

# New column G ("Reference") values, filled row-by-row (top to bottom) to
# mirror the order the references were typed in during the original edit.
$refs = @{
    3  = "Anker-Nilssen and Aavark 2006"
    4  = $ram
    5  = $ram
    6  = $ram
    7  = $ram
    8  = $ram
    9  = $ram
    10 = $ram
    11 = $ram
    12 = $ram
    13 = $ram
    14 = $ram
    15 = $ram
    16 = $ram
    17 = $ram
    18 = $ram
    19 = $ram
    20 = $ram
    21 = $ram
    22 = $ram
    23 = $ram
    24 = $ram
    25 = "ICES 2016"
    26 = "Furness 2007"
    27 = $ram
    28 = $ram
    29 = $ram
    30 = $ram
    32 = $ram
    33 = $ram
    34 = $ram
    35 = $ram
    37 = "Crawford et al. 2006"
    38 = "Crawford et al. 2006"
    39 = "Crawford et al. 2006"
    40 = "Crawford et al. 2006"
    41 = $ram
    43 = $ram
    44 = $ram
    45 = $ram
    46 = $ram
    47 = $ram
    48 = $ram
    49 = "NEFSC 2018"
    50 = $ram
    51 = $ram
    52 = $ram
    53 = $ram
    54 = "Hendrickson and Showell 2016"
    55 = $ram
    57 = "MacCall et al. 2016"
    58 = $ram
    59 = $ram
    60 = $ram
    61 = "Mills et al. 2007"
    62 = "Mills et al. 2007"
    63 = "Mills et al. 2007"
    64 = $ram
    65 = $ram
    66 = $ram
    67 = $ram
}

for ($r = 3; $r -le 67; $r++) {
    if ($refs.ContainsKey($r)) {
        $ws.Cells.Item($r, 7).Value = $refs[$r]
    }
}

# New column header, added last (matches shared-string append order).
# Copy F1's header formatting (bold font + thick-bottom border) onto G1.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Reference"
$excel.CutCopyMode = $false

# Column G width to fit the new content (stored width in the xlsx is
# ColumnWidth + ~0.8333, so back that padding out to land on width="27")
$ws.Columns.Item(7).ColumnWidth = 26.1666666666667

# View state change captured in the diff: move the active selection to G62
# (the top-left scroll anchor resets to default as a side effect)
$ws.Range("G62").Select()
